# Auto-generated edit script: update Leve profit-calculation data per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 256.66666
$ws.Range("I6").Value = 12.5
$ws.Range("J6").Value = 378.75
$ws.Range("K6").Value = 37.5
$ws.Range("L6").Value = 1136.25
$ws.Range("M6").Value = 74.5
$ws.Range("N6").Value = -1360.25

$ws.Range("H8").Value = 121
$ws.Range("I8").Value = 123.625
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 370.875
$ws.Range("L8").Value = 300
$ws.Range("M8").Value = -231.875
$ws.Range("N8").Value = -578

$ws.Range("H69").Value = 1050
$ws.Range("I69").Value = 1050
$ws.Range("K69").Value = 3150
$ws.Range("M69").Value = -2276

$ws.Range("H72").Value = 1050
$ws.Range("I72").Value = 1050
$ws.Range("K72").Value = 9450
$ws.Range("M72").Value = -5082

$ws.Range("H98").Value = 2000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 2000
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -4996

$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -10900

$ws.Range("H132").Value = 6233
$ws.Range("I132").Value = 6370.1113
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 19110.3339
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -16580.3339
$ws.Range("N132").Value = -20057

$ws.Range("H137").Value = 2797.8
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 3994.5
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 11983.5
$ws.Range("M137").Value = -3450
$ws.Range("N137").Value = -17083.5

$ws.Range("H141").Value = 18526.715
$ws.Range("I141").Value = 28671.75
$ws.Range("K141").Value = 86015.25
$ws.Range("M141").Value = -80835.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20700

$ws.Range("H17").Value = 2250
$ws.Range("I17").Value = 500
$ws.Range("K17").Value = 500
$ws.Range("M17").Value = -327

$ws.Range("H42").Value = 10028
$ws.Range("I42").Value = 10028
$ws.Range("K42").Value = 10028
$ws.Range("M42").Value = -9542

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H61").Value = 42000
$ws.Range("I61").Value = 42000
$ws.Range("K61").Value = 42000
$ws.Range("M61").Value = -41788

$ws.Range("H74").Value = 1524
$ws.Range("I74").Value = 1033.8334
$ws.Range("K74").Value = 1033.8334
$ws.Range("M74").Value = -159.8334

$ws.Range("H77").Value = 1524
$ws.Range("I77").Value = 1033.8334
$ws.Range("K77").Value = 5169.166999999999
$ws.Range("M77").Value = -801.1669999999995

$ws.Range("H92").Value = 53920
$ws.Range("J92").Value = 53920
$ws.Range("L92").Value = 53920
$ws.Range("N92").Value = -58912

$ws.Range("H136").Value = 42000
$ws.Range("I136").Value = 42000
$ws.Range("K136").Value = 126000
$ws.Range("M136").Value = -123450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 272.5
$ws.Range("I10").Value = 350
$ws.Range("J10").Value = 195
$ws.Range("K10").Value = 350
$ws.Range("L10").Value = 195
$ws.Range("M10").Value = -210
$ws.Range("N10").Value = -475

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H134").Value = 5609
$ws.Range("I134").Value = 2010.3334
$ws.Range("K134").Value = 6031.0002
$ws.Range("M134").Value = -3496.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 316.66666
$ws.Range("I10").Value = 316.66666
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 316.66666
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -177.66666
$ws.Range("N10").ClearContents()

$ws.Range("H22").Value = 1153.3334
$ws.Range("I22").Value = 729.25
$ws.Range("J22").Value = 2001.5
$ws.Range("K22").Value = 729.25
$ws.Range("L22").Value = 2001.5
$ws.Range("M22").Value = -379.25
$ws.Range("N22").Value = -2701.5

$ws.Range("H31").Value = 3423.5
$ws.Range("I31").Value = 2150
$ws.Range("K31").Value = 2150
$ws.Range("M31").Value = -1855

$ws.Range("H34").Value = 3423.5
$ws.Range("I34").Value = 2150
$ws.Range("K34").Value = 2150
$ws.Range("M34").Value = -1948

$ws.Range("H47").Value = 22150
$ws.Range("I47").Value = 19300
$ws.Range("K47").Value = 19300
$ws.Range("M47").Value = -18734

$ws.Range("H57").Value = 45000
$ws.Range("J57").Value = 45000
$ws.Range("L57").Value = 45000
$ws.Range("N57").Value = -46120

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H99").Value = 2719.5715
$ws.Range("I99").Value = 1717.4
$ws.Range("K99").Value = 1717.4
$ws.Range("M99").Value = -219.4000000000001

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").ClearContents()

$ws.Range("H126").Value = 2719.5715
$ws.Range("I126").Value = 1717.4
$ws.Range("K126").Value = 5152.200000000001
$ws.Range("M126").Value = -2682.200000000001

$ws.Range("H132").Value = 6236.625
$ws.Range("I132").Value = 4483.1665
$ws.Range("J132").Value = 11497
$ws.Range("K132").Value = 13449.4995
$ws.Range("L132").Value = 34491
$ws.Range("M132").Value = -10919.4995
$ws.Range("N132").Value = -39551

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 159.6
$ws.Range("J7").Value = 174.75
$ws.Range("L7").Value = 524.25
$ws.Range("N7").Value = -748.25

$ws.Range("H119").Value = 3104.6667
$ws.Range("I119").Value = 3104.6667
$ws.Range("K119").Value = 9314.000100000001
$ws.Range("M119").Value = -4476.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 19333.334
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 27500
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 27500
$ws.Range("M46").Value = -2844
$ws.Range("N46").Value = -27812

$ws.Range("H70").Value = 6999.25
$ws.Range("I70").Value = 6999
$ws.Range("J70").Value = 6999.5
$ws.Range("K70").Value = 6999
$ws.Range("L70").Value = 6999.5
$ws.Range("M70").Value = -6729
$ws.Range("N70").Value = -7539.5

$ws.Range("H73").Value = 6999.25
$ws.Range("I73").Value = 6999
$ws.Range("J73").Value = 6999.5
$ws.Range("K73").Value = 6999
$ws.Range("L73").Value = 6999.5
$ws.Range("M73").Value = -6063
$ws.Range("N73").Value = -8871.5

$ws.Range("H80").Value = 20602.875
$ws.Range("I80").Value = 18103.25
$ws.Range("K80").Value = 18103.25
$ws.Range("M80").Value = -17105.25

$ws.Range("H83").Value = 20602.875
$ws.Range("I83").Value = 18103.25
$ws.Range("K83").Value = 90516.25
$ws.Range("M83").Value = -85524.25

$ws.Range("H132").Value = 1800
$ws.Range("I132").Value = 1800
$ws.Range("K132").Value = 5400
$ws.Range("M132").Value = -2870

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1968
$ws.Range("J61").Value = 1200
$ws.Range("L61").Value = 1200
$ws.Range("N61").Value = -1604

$ws.Range("H113").Value = 1968
$ws.Range("J113").Value = 1200
$ws.Range("L113").Value = 1200
$ws.Range("N113").Value = -5540

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 25000
$ws.Range("I70").Value = 25000
$ws.Range("K70").Value = 25000
$ws.Range("M70").Value = -24685

$ws.Range("H73").Value = 25000
$ws.Range("I73").Value = 25000
$ws.Range("K73").Value = 25000
$ws.Range("M73").Value = -23908
